$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rotates the per-observation data (Id, Antal, Enhet, Alder-Stadium,
# Ost, Nord, Publik kommentar) among rows 2, 3, 4 and 6 in a 4-cycle:
#   new row2 <- old row6
#   new row3 <- old row2
#   new row4 <- old row3
#   new row6 <- old row4
# Row 5 (and every other column) is left untouched.

function Set-RowData($r, $data) {
    $ws.Range("A$r").Value = $data.A

    # Columns I/J/K hold numeric-looking values stored as text in the source
    # file (t="inlineStr"); force text storage so "7", "30", etc. don't get
    # reinterpreted as numbers.
    $ws.Range("I$r").NumberFormat = "@"
    $ws.Range("I$r").Value = $data.I

    $ws.Range("J$r").Value = $data.J

    if ($null -eq $data.K) {
        $ws.Range("K$r").ClearContents()
    } else {
        $ws.Range("K$r").Value = $data.K
    }

    $ws.Range("Q$r").Value = $data.Q
    $ws.Range("R$r").Value = $data.R

    if ($null -eq $data.AC) {
        $ws.Range("AC$r").ClearContents()
    } else {
        $ws.Range("AC$r").Value = $data.AC
    }
}

# Original (pre-edit) contents of the four affected rows, literal per the
# workbook's source data.
$row2 = @{ A = 111416528; I = "30"; J = "plantor/tuvor"; K = $null; Q = 359092.1819271583; R = 6393204.710604292; AC = $null }
$row3 = @{ A = 111416525; I = "4"; J = "stjälkar/strån/skott"; K = "blomning"; Q = 359095.1406046218; R = 6393212.639220579; AC = "även bladrosetter på 1 kvm" }
$row4 = @{ A = 111416521; I = "1"; J = "stjälkar/strån/skott"; K = "blomning"; Q = 359101.3469427949; R = 6393205.997596246; AC = "även ca 30 bladrosetter" }
$row6 = @{ A = 111416523; I = "7"; J = "stjälkar/strån/skott"; K = "blomning"; Q = 359100.0376043977; R = 6393214.610374114; AC = "även bladrosetter på 1 kvm" }

Set-RowData 2 $row6
Set-RowData 3 $row2
Set-RowData 4 $row3
Set-RowData 6 $row4
